# Did you know #4 — refresh the "captured" date from 5/9/2018 to 5/14/2018
# across the handout master, notes master and all title-slide layouts, and
# update the slide's call-to-action subtitle text.

$p = $ppt.ActivePresentation

# --- 1. Main content edit: slide 1 subtitle ------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "#4 Using a workspace") {
        $shp.TextFrame.TextRange.Text = "#4 Learn more…"
    }
}

# --- 2. Auto-update date placeholders (Handout Master / Notes Master) ----
# These are genuine "dt" placeholders wired through HeadersFooters.DateAndTime
# (the cached field text that was baked in on 5/9/2018, now refreshed).
$oldDateShort = "5/9/2018"
$newDateShort = "5/14/2018"

$handoutMaster = $p.HandoutMaster
$handoutMaster.HeadersFooters.DateAndTime.Text = $newDateShort

$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = $newDateShort

# --- 3. "May 9, 2018" style date text boxes on the title-slide layouts ---
# Every CustomLayout that carries the small printed-date text box in the
# corner gets its cached text bumped the same five days forward.
$oldDateLong = "May 9, 2018"
$newDateLong = "May 14, 2018"

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq $oldDateLong) {
            $shp.TextFrame.TextRange.Text = $newDateLong
        }
    }
}
